$d = $word.ActiveDocument

$find = "Perioadele campaniei din Constelația Gemeni 2022: 14-23 februarie, 14-24 martie"
$replace = "Perioadele campaniei din 2022 pentru Constelația Gemeni: 14-23 februarie, 14-24 martie"

$r = $d.Content
$r.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
